$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.564.01"
$ws.Range("E2").Value = "  -0.78%  "

# Row 3
$ws.Range("D3").Value = "3.545.22"
$ws.Range("E3").Value = "  -2.13%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "585.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.12%  "

# Row 7
$ws.Range("E7").Value = "  -2.01%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("E9").Value = "  +0.12%  "

# Row 10
$ws.Range("E10").Value = "  -3.26%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.18"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.32%  "

# Row 12
$ws.Range("E12").Value = "  -5.06%  "

# Row 13
$ws.Range("B13").Value = "BitcoinCash"
$ws.Range("C13").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "688.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +16.50%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.37%  "

# Row 15
$ws.Range("D15").Value = "4.108.84"
$ws.Range("E15").Value = "  -2.00%  "

# Row 16
$ws.Range("D16").Value = "69.655.53"
$ws.Range("E16").Value = "  -0.83%  "

# Row 17
$ws.Range("D17").Value = "3.545.83"
$ws.Range("E17").Value = "  -2.06%  "

# Row 18
$ws.Range("E18").Value = "  -5.97%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.62"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.24%  "

# Row 20
$ws.Range("E20").Value = "  -0.71%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.973"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.22%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.91%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "108.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +5.80%  "

# Row 24
$ws.Range("E24").Value = "  +0.97%  "

# Row 25
$ws.Range("E25").Value = "  -4.62%  "

# Row 26
$ws.Range("E26").Value = "  -2.90%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.28"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.50%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.25%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.55%  "

# Row 30
$ws.Range("E30").Value = "  -1.66%  "

# Row 31
$ws.Range("E31").Value = "  -3.41%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.16%  "

# Row 33
$ws.Range("E33").Value = "  -4.24%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "62.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.82%  "

# Row 35
$ws.Range("D35").Value = "3.809.99"
$ws.Range("E35").Value = "  -3.74%  "

# Row 36
$ws.Range("E36").Value = "  -8.85%  "

# Row 37
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.70"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.42%  "

# Row 38
$ws.Range("B38").Value = "Dai"
$ws.Range("C38").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.12%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.54%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "499.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.38%  "

# Row 41
$ws.Range("E41").Value = "  -4.66%  "

# Row 42
$ws.Range("E42").Value = "  +1.26%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "34.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.65%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0460"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.88%  "

# Row 45
$ws.Range("E45").Value = "  +2.71%  "

# Row 46
$ws.Range("E46").Value = "  +1.09%  "

# Row 47
$ws.Range("E47").Value = "  -2.31%  "

# Row 48
$ws.Range("E48").Value = "  -0.26%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.39%  "

# Row 50
$ws.Range("E50").Value = "  +21.79%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +67.67%  "

